# Apply "changed assignment file and added IDE and PL" edit:
# Adds two new columns (E = Deadline, F = Отметка о выполении / mark of
# completion) with date values + number formatting to the assignment
# tracking sheet, and a "-" marker in E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# Column headers (row 1).
# F1 is written before E1 so that the new shared strings end up in the
# same order as the reference edit (Отметка о выполении, Deadline, -).
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Отметка о выполении"
$ws.Range("F1").NumberFormat = "dd/mm/yy\ h:mm;@"
$ws.Range("F1").HorizontalAlignment = $xlCenter

$ws.Range("E1").Value = "Deadline"
$ws.Range("D3").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Column E - Deadline date/time values (rows 2-15), centered, built-in
# date-time number format.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = 42273.833333333336
$ws.Range("E2").NumberFormat = "m/d/yy h:mm"
$ws.Range("E2").HorizontalAlignment = $xlCenter
$ws.Range("E2").VerticalAlignment = $xlCenter

$ws.Range("E2").Copy()

$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 42273.833333333336

$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 42273.833333333336

$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 42273.833333333336

$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 42273.833333333336

$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 42273.833333333336

$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 42273.833333333336

$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = 42273.833333333336

$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = 42273.833333333336

$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = 42273.833333333336

$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = 42273.833333333336

$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = 42273.833333333336

$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = 42274.958333333336

$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 42270.833333333336

# E16 - text marker "-" (reuses the plain center/center alignment style).
$ws.Range("E16").Value = "-"
$ws.Range("D3").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Column F - "mark of completion" date/time values, centered (custom
# dd/mm/yy h:mm format), only present for rows 11, 12 and 15.
# ---------------------------------------------------------------------
$ws.Range("F11").Value = 42272.979166666664
$ws.Range("F1").Copy()
$ws.Range("F11").PasteSpecial(-4122)

$ws.Range("F12").Value = 42272.979166666664
$ws.Range("F1").Copy()
$ws.Range("F12").PasteSpecial(-4122)

$ws.Range("F15").Value = 42269.083333333336
$ws.Range("F1").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Column widths for the two new columns.
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 19.8
$ws.Columns("F").ColumnWidth = 20.6

# ---------------------------------------------------------------------
# Selection / view state.
# ---------------------------------------------------------------------
$ws.Range("D16").Select()
